# update week 71 uk
# Adds the two new survey rows (wave 26, survey_round 70 & 71) that were
# received for the UK panel, following the same pattern as the existing
# rows in the "UK" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

# --- Row 83: survey_round 70, panel E, wave 26 ---------------------------
$ws.Range("A83").Value = 3
$ws.Range("B83").Value = 0
$ws.Range("C83").Value = "uk"
$ws.Range("D83").Value = 70
$ws.Range("E83").Value = "E"
$ws.Range("F83").Value = 26

# Copy the date formatting from the row above (so the new cell keeps the
# existing "date_recieved" number format instead of creating a new one),
# then set the actual date value.
$ws.Range("G82").Copy()
$ws.Range("G83").PasteSpecial(-4122)
$ws.Range("G83").Value = "2021-07-30"

$ws.Range("H83").Value = "21-037558_PEW26_Final_ICUO"

# --- Row 84: survey_round 71, panel F, wave 26 ---------------------------
$ws.Range("A84").Value = 3
$ws.Range("B84").Value = 0
$ws.Range("C84").Value = "uk"
$ws.Range("D84").Value = 71
$ws.Range("E84").Value = "F"
$ws.Range("F84").Value = 26

$ws.Range("G82").Copy()
$ws.Range("G84").PasteSpecial(-4122)
$ws.Range("G84").Value = "2021-08-05"

$ws.Range("H84").Value = "21-037554_PFW26_Final_ICUO"

# --- Column I: same "r_name" formula pattern used by every other row -----
$ws.Range("I83:I84").Formula = "=C83&""_""&""sr""&TEXT(D83,""00"")&""_""&YEAR(G83)&TEXT(G83,""MM"")&TEXT(G83,""DD"")&""_p""&E83&""_wv""&TEXT(F83,""00"")&"""""

# Move the selection the way the original author left it after entering
# the new data.
$ws.Range("I88").Select()
